# Auto-generated edit script: updates Phantom_Profits price/profit figures
# across all 8 item-category sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 281.8889
$ws.Range("I2").Value = 288.2857
$ws.Range("K2").Value = 288.2857
$ws.Range("M2").Value = -175.2857
$ws.Range("H17").Value = 2998
$ws.Range("J17").Value = 2998
$ws.Range("L17").Value = 8994
$ws.Range("N17").Value = -9330
$ws.Range("H69").Value = 15626.667
$ws.Range("I69").Value = 10500
$ws.Range("J69").Value = 16908.334
$ws.Range("K69").Value = 31500
$ws.Range("L69").Value = 50725.00199999999
$ws.Range("M69").Value = -30626
$ws.Range("N69").Value = -52473.00199999999
$ws.Range("H72").Value = 15626.667
$ws.Range("I72").Value = 10500
$ws.Range("J72").Value = 16908.334
$ws.Range("K72").Value = 94500
$ws.Range("L72").Value = 152175.006
$ws.Range("M72").Value = -90132
$ws.Range("N72").Value = -160911.006
$ws.Range("H100").Value = 4548.625
$ws.Range("I100").Value = 4548.625
$ws.Range("K100").Value = 4548.625
$ws.Range("M100").Value = -4007.625
$ws.Range("H111").Value = 3644.111
$ws.Range("I111").Value = 4282.8335
$ws.Range("J111").Value = 2366.6667
$ws.Range("K111").Value = 12848.5005
$ws.Range("L111").Value = 7100.000100000001
$ws.Range("M111").Value = -9781.500499999998
$ws.Range("N111").Value = -13234.0001
$ws.Range("H116").Value = 5096.2
$ws.Range("J116").Value = 4496.5
$ws.Range("L116").Value = 4496.5
$ws.Range("N116").Value = -11380.5
$ws.Range("H137").Value = 25642138
$ws.Range("I137").Value = 27778816
$ws.Range("K137").Value = 83336448
$ws.Range("M137").Value = -83333898

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3720.1738
$ws.Range("I32").Value = 2646.1428
$ws.Range("J32").Value = 14997.5
$ws.Range("K32").Value = 2646.1428
$ws.Range("L32").Value = 14997.5
$ws.Range("M32").Value = -2359.1428
$ws.Range("N32").Value = -15571.5
$ws.Range("H61").Value = 1514.7391
$ws.Range("I61").Value = 1538.8422
$ws.Range("K61").Value = 1538.8422
$ws.Range("M61").Value = -1326.8422
$ws.Range("H110").Value = 498
$ws.Range("I110").Value = 498
$ws.Range("K110").Value = 498
$ws.Range("M110").Value = 1547
$ws.Range("H132").Value = 3763.9546
$ws.Range("I132").Value = 2721.3572
$ws.Range("K132").Value = 8164.071599999999
$ws.Range("M132").Value = -5634.071599999999
$ws.Range("H136").Value = 1514.7391
$ws.Range("I136").Value = 1538.8422
$ws.Range("K136").Value = 4616.5266
$ws.Range("M136").Value = -2066.5266

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3373.5557
$ws.Range("I105").Value = 3370.25
$ws.Range("K105").Value = 3370.25
$ws.Range("M105").Value = -1623.25
$ws.Range("H134").Value = 1634.6
$ws.Range("I134").Value = 1624.826
$ws.Range("K134").Value = 4874.478
$ws.Range("M134").Value = -2339.478

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 823.2727
$ws.Range("I16").Value = 823.2727
$ws.Range("K16").Value = 823.2727
$ws.Range("M16").Value = -536.2727
$ws.Range("H31").Value = 6153.1904
$ws.Range("I31").Value = 10347.125
$ws.Range("J31").Value = 3572.3076
$ws.Range("K31").Value = 10347.125
$ws.Range("L31").Value = 3572.3076
$ws.Range("M31").Value = -10052.125
$ws.Range("N31").Value = -4162.3076
$ws.Range("H34").Value = 6153.1904
$ws.Range("I34").Value = 10347.125
$ws.Range("J34").Value = 3572.3076
$ws.Range("K34").Value = 10347.125
$ws.Range("L34").Value = 3572.3076
$ws.Range("M34").Value = -10145.125
$ws.Range("N34").Value = -3976.3076
$ws.Range("H62").Value = 2790.5
$ws.Range("J62").Value = 2309.6667
$ws.Range("L62").Value = 2309.6667
$ws.Range("N62").Value = -3557.6667
$ws.Range("H65").Value = 2790.5
$ws.Range("J65").Value = 2309.6667
$ws.Range("L65").Value = 11548.3335
$ws.Range("N65").Value = -17788.3335
$ws.Range("H74").Value = 69999
$ws.Range("J74").Value = 69999
$ws.Range("L74").Value = 69999
$ws.Range("N74").Value = -71747
$ws.Range("H77").Value = 69999
$ws.Range("J77").Value = 69999
$ws.Range("L77").Value = 209997
$ws.Range("N77").Value = -218733
$ws.Range("H107").Value = 1324.2941
$ws.Range("I107").Value = 1116.6154
$ws.Range("K107").Value = 1116.6154
$ws.Range("M107").Value = 803.3846000000001
$ws.Range("H113").Value = 823.2727
$ws.Range("I113").Value = 823.2727
$ws.Range("K113").Value = 823.2727
$ws.Range("M113").Value = 1346.7273
$ws.Range("H122").Value = 1340.75
$ws.Range("I122").Value = 1354.0714
$ws.Range("K122").Value = 4062.2142
$ws.Range("M122").Value = -1612.2142
$ws.Range("H132").Value = 28576910
$ws.Range("I132").Value = 50004900
$ws.Range("K132").Value = 150014700
$ws.Range("M132").Value = -150012170

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1894.3771
$ws.Range("J131").Value = 1883.807
$ws.Range("L131").Value = 5651.421
$ws.Range("N131").Value = -15731.421

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5213.778
$ws.Range("I70").Value = 5326.4287
$ws.Range("J70").Value = 4819.5
$ws.Range("K70").Value = 5326.4287
$ws.Range("L70").Value = 4819.5
$ws.Range("M70").Value = -5056.4287
$ws.Range("N70").Value = -5359.5
$ws.Range("H73").Value = 5213.778
$ws.Range("I73").Value = 5326.4287
$ws.Range("J73").Value = 4819.5
$ws.Range("K73").Value = 5326.4287
$ws.Range("L73").Value = 4819.5
$ws.Range("M73").Value = -4390.4287
$ws.Range("N73").Value = -6691.5
$ws.Range("H96").Value = 19340.25
$ws.Range("I96").Value = 4000
$ws.Range("J96").Value = 24453.666
$ws.Range("K96").Value = 4000
$ws.Range("L96").Value = 24453.666
$ws.Range("M96").Value = -1254
$ws.Range("N96").Value = -29945.666
$ws.Range("H113").Value = 6732.2
$ws.Range("I113").Value = 5665.25
$ws.Range("K113").Value = 5665.25
$ws.Range("M113").Value = -3495.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4374.3335
$ws.Range("I46").Value = 3220
$ws.Range("J46").Value = 5143.8887
$ws.Range("K46").Value = 3220
$ws.Range("L46").Value = 5143.8887
$ws.Range("M46").Value = -3032
$ws.Range("N46").Value = -5519.8887
$ws.Range("H68").Value = 4388
$ws.Range("I68").Value = 1060.25
$ws.Range("K68").Value = 1060.25
$ws.Range("M68").Value = -311.25
$ws.Range("H69").Value = 11998
$ws.Range("J69").Value = 11998
$ws.Range("L69").Value = 11998
$ws.Range("N69").Value = -13620
$ws.Range("H71").Value = 4388
$ws.Range("I71").Value = 1060.25
$ws.Range("K71").Value = 5301.25
$ws.Range("M71").Value = -1557.25
$ws.Range("H72").Value = 11998
$ws.Range("J72").Value = 11998
$ws.Range("L72").Value = 35994
$ws.Range("N72").Value = -44106
$ws.Range("H93").Value = 531.8889
$ws.Range("I93").Value = 600.25
$ws.Range("J93").Value = 520
$ws.Range("K93").Value = 600.25
$ws.Range("L93").Value = 520
$ws.Range("M93").Value = 647.75
$ws.Range("N93").Value = -3016
$ws.Range("H122").Value = 2927.4443
$ws.Range("I122").Value = 2658.1667
$ws.Range("K122").Value = 7974.500100000001
$ws.Range("M122").Value = -5524.500100000001
$ws.Range("H132").Value = 1999.75
$ws.Range("I132").Value = 999.6667
$ws.Range("K132").Value = 2999.0001
$ws.Range("M132").Value = -469.0001000000002

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H113").Value = 390.85715
$ws.Range("I113").Value = 315.75
$ws.Range("J113").Value = 491
$ws.Range("K113").Value = 947.25
$ws.Range("L113").Value = 1473
$ws.Range("M113").Value = 1222.75
$ws.Range("N113").Value = -5813
$ws.Range("H132").Value = 333337950
$ws.Range("I132").Value = 6943.5
$ws.Range("J132").Value = 1000000000
$ws.Range("K132").Value = 20830.5
$ws.Range("L132").Value = 3000000000
$ws.Range("M132").Value = -18300.5
$ws.Range("N132").Value = -3000005060
$ws.Range("H140").Value = 54758.855
$ws.Range("J140").Value = 54758.855
$ws.Range("L140").Value = 54758.855
$ws.Range("N140").Value = -65118.855
